$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '45.724.17'
$ws.Range('E2').Value = '  +2.62%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.436.51'
$ws.Range('E3').Value = '  +0.20%  '
$ws.Range('E4').Value = '  -0.06%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '321.61'
$ws.Range('E5').Value = '  +3.89%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '104.32'
$ws.Range('E6').Value = '  +2.29%  '
$ws.Range('E7').Value = '  +0.84%  '
$ws.Range('E8').Value = '  -0.01%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.535'
$ws.Range('E9').Value = '  +5.41%  '
$ws.Range('E10').Value = '  +1.51%  '
$ws.Range('E11').Value = '  +0.49%  '
$ws.Range('E12').Value = '  -1.27%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '18.26'
$ws.Range('E13').Value = '  -2.86%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '7.03'
$ws.Range('E14').Value = '  +0.78%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '2.816.74'
$ws.Range('E15').Value = '  +0.14%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '2.434.35'
$ws.Range('E16').Value = '  +0.42%  '
$ws.Range('E17').Value = '  +0.06%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '45.603.77'
$ws.Range('E18').Value = '  +2.54%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '12.41'
$ws.Range('E19').Value = '  -0.41%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '6.42'
$ws.Range('E20').Value = '  +0.19%  '
$ws.Range('E21').Value = '  +2.26%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '71.69'
$ws.Range('E22').Value = '  +4.07%  '
$ws.Range('E23').Value = '  +1.34%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.50'
$ws.Range('E25').Value = '  +0.29%  '
$ws.Range('E26').Value = '  +0.02%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '25.75'
$ws.Range('E27').Value = '  +2.17%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.28'
$ws.Range('E28').Value = '  -0.50%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '9.66'
$ws.Range('E29').Value = '  -0.24%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '33.49'
$ws.Range('E30').Value = '  +0.88%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '49.28'
$ws.Range('E31').Value = '  +1.58%  '
$ws.Range('E32').Value = '  +4.87%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '20.14'
$ws.Range('E33').Value = '  +2.57%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '5.25'
$ws.Range('E34').Value = '  +1.02%  '
$ws.Range('E35').Value = '  +0.13%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.0759'
$ws.Range('E36').Value = '  -0.72%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '4.55'
$ws.Range('E37').Value = '  +1.10%  '
$ws.Range('E38').Value = '  -0.80%  '
$ws.Range('E39').Value = '  +0.57%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '127.42'
$ws.Range('E40').Value = '  +0.92%  '
$ws.Range('E41').Value = '  -2.41%  '
$ws.Range('E42').Value = '  +1.32%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '20.90'
$ws.Range('E43').Value = '  -4.95%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.0291'
$ws.Range('E44').Value = '  +0.18%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '1.959.39'
$ws.Range('E45').Value = '  +0.58%  '
$ws.Range('E46').Value = '  +0.93%  '
$ws.Range('E47').Value = '  -2.79%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.81'
$ws.Range('E48').Value = '  +8.00%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '9.13'
$ws.Range('E49').Value = '  -4.39%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '77.27'
$ws.Range('E50').Value = '  +4.41%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '4.85'
$ws.Range('E51').Value = '  +4.83%  '
